# Add "Epic" and "Storyboard- Backlog" worksheets after the existing
# "Storyboard- Active Sprint" sheet, populate them with the Field
# Specification rows, and make "Storyboard- Backlog" the active tab.

$wb = $excel.ActiveWorkbook
$active = $wb.Worksheets.Item(1)

# ---------------------------------------------------------------------
# Create the new sheets, in order, right after the active sprint sheet.
# ---------------------------------------------------------------------
$epic = $wb.Worksheets.Add($null, $active)
$epic.Name = "Epic"

$backlog = $wb.Worksheets.Add($null, $epic)
$backlog.Name = "Storyboard- Backlog"

# ---------------------------------------------------------------------
# Epic sheet content
# ---------------------------------------------------------------------
$epic.Range("A1").Value = "Field Name"
$epic.Range("B1").Value = "Type"
$epic.Range("C1").Value = "Spec ID"
$epic.Range("D1").Value = "Label"
$epic.Range("E1").Value = "Placeholder/Default Selection"
$epic.Range("F1").Value = "Description"
$epic.Range("G1").Value = "Data Set"
$epic.Range("H1").Value = "Allowed Values"
$epic.Range("I1").Value = "Error Scenario"
$epic.Range("J1").Value = "Error Message"
$epic.Range("K1").Value = "Error code"

$epic.Range("H2").WrapText = $true

$epic.Range("A3").Value = "sortBy"
$epic.Range("B3").Value = "Select"
$epic.Range("C3").Value = "EM_E_1"
$epic.Range("D3").Value = "N/A"
$epic.Range("E3").Value = "Sort By"
$epic.Range("F3").Value = "Options for sorting Epics on the Epics page"
$epic.Range("G3").Value = "N/A"
$epic.Range("H3").Value = "Product`nOwner"
$epic.Range("H3").WrapText = $true
$epic.Range("I3").Value = "N/A"
$epic.Range("J3").Value = "N/A"
$epic.Range("K3").Value = "N/A"
$epic.Rows.Item(3).RowHeight = 30

$epic.Range("A5").Value = "filter"
$epic.Range("B5").Value = "Textfield"
$epic.Range("C5").Value = "EM_E_2"
$epic.Range("D5").Value = "N/A"
$epic.Range("E5").Value = "Filter"
$epic.Range("F5").Value = "textfield to enter text to filter Epics"
$epic.Range("G5").Value = "Alphanumeric`nSpecial chars : Spaces "
$epic.Range("G5").WrapText = $true
$epic.Range("H5").Value = "N/A"
$epic.Range("I5").Value = "Special characters except for spaces are entered"
$epic.Range("J5").Value = "Enter a valid value"
$epic.Range("K5").Value = "ES_008"
$epic.Rows.Item(5).RowHeight = 60

$epic.Columns.Item(5).ColumnWidth = 24.16
$epic.Columns.Item(6).ColumnWidth = 34.66
$epic.Columns.Item(7).ColumnWidth = 11.5
$epic.Columns.Item(8).ColumnWidth = 13.83
$epic.Columns.Item(9).ColumnWidth = 37.66

$epic.Activate() | Out-Null
$epic.Rows.Item(1).Select() | Out-Null

# ---------------------------------------------------------------------
# Storyboard- Backlog sheet content
# ---------------------------------------------------------------------
$backlog.Range("A1").Value = "Field Name"
$backlog.Range("B1").Value = "Type"
$backlog.Range("C1").Value = "Spec ID"
$backlog.Range("D1").Value = "Label"
$backlog.Range("E1").Value = "Placeholder/Default Selection"
$backlog.Range("F1").Value = "Description"
$backlog.Range("G1").Value = "Data Set"
$backlog.Range("H1").Value = "Allowed Values"
$backlog.Range("I1").Value = "Error Scenario"
$backlog.Range("J1").Value = "Error Message"
$backlog.Range("K1").Value = "Error code"

$backlog.Range("A2").Value = "Project"
$backlog.Range("B2").Value = "Select"
$backlog.Range("C2").Value = "SB_B_1"
$backlog.Range("D2").Value = "N/A"
$backlog.Range("E2").Value = "Project"
$backlog.Range("F2").Value = "Filter issues by selecting a project"
$backlog.Range("G2").Value = "N/A"
$backlog.Range("H2").Value = "Project names"
$backlog.Range("H2").WrapText = $true
$backlog.Range("I2").Value = "N/A"
$backlog.Range("J2").Value = "N/A"
$backlog.Range("K2").Value = "N/A"

$backlog.Range("A3").Value = "filterByText"
$backlog.Range("B3").Value = "Textfield"
$backlog.Range("C3").Value = "SB_B_2"
$backlog.Range("D3").Value = "N/A"
$backlog.Range("E3").Value = "Contain text"
$backlog.Range("F3").Value = "textfield to enter text to filter issues"
$backlog.Range("G3").Value = "Alphanumeric`nSpecial chars : Spaces "
$backlog.Range("G3").WrapText = $true
$backlog.Range("H3").Value = "N/A"
$backlog.Range("I3").Value = "Special characters except for spaces are entered"
$backlog.Range("J3").Value = "Enter a valid value"
$backlog.Range("K3").Value = "ES_008"
$backlog.Rows.Item(3).RowHeight = 45

$backlog.Range("A4").Value = "Type/Owner/Status"
$backlog.Range("B4").Value = "Select in 3 buttons"
$backlog.Range("C4").Value = "SB_B_3"
$backlog.Range("D4").Value = "N/A"
$backlog.Range("E4").Value = "Type"
$backlog.Range("F4").Value = "Options on navigation bar for sorting issues on the backlog page"
$backlog.Range("G4").Value = "N/A"
$backlog.Range("H4").Value = "Type`nOwner`nStatus"
$backlog.Range("H4").WrapText = $true
$backlog.Range("I4").Value = "N/A"
$backlog.Range("J4").Value = "N/A"
$backlog.Range("K4").Value = "N/A"
$backlog.Rows.Item(4).RowHeight = 45

$backlog.Columns.Item(1).ColumnWidth = 16.16
$backlog.Columns.Item(2).ColumnWidth = 15.66
$backlog.Columns.Item(3).ColumnWidth = 10.33
$backlog.Columns.Item(4).ColumnWidth = 5.5
$backlog.Columns.Item(5).ColumnWidth = 25.5
$backlog.Columns.Item(6).ColumnWidth = 50.16
$backlog.Columns.Item(7).ColumnWidth = 12.66
$backlog.Columns.Item(8).ColumnWidth = 13.66
$backlog.Columns.Item(9).ColumnWidth = 48.5
$backlog.Columns.Item(10).ColumnWidth = 16.16
$backlog.Columns.Item(11).ColumnWidth = 9.83

$backlog.PageSetup.Orientation = 1

$backlog.Activate() | Out-Null
$backlog.Range("B9").Select() | Out-Null
